# Add two new columns, I ("I0") and J ("IF"), to the sheet, mirroring the
# style of the existing header row and filling in the per-row numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells: set the new header labels, then copy the style used by the
# other header cells (e.g. H1) onto them.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2-44: numeric values for columns I (I0) and J (IF).
$data = @(
    @(2, 7, 7),
    @(3, 7, 7),
    @(4, 7, 7),
    @(5, 8, 8),
    @(6, 8, 9),
    @(7, 8, 8),
    @(8, 8, 8),
    @(9, 7, 8),
    @(10, 7, 7),
    @(11, 7, 7),
    @(12, 7, 8),
    @(13, 7, 7),
    @(14, 8, 8),
    @(15, 7, 7),
    @(16, 8, 8),
    @(17, 8, 8),
    @(18, 8, 9),
    @(19, 7, 7),
    @(20, 8, 8),
    @(21, 8, 8),
    @(22, 8, 8),
    @(23, 9, 9),
    @(24, 7, 7),
    @(25, 6, 6),
    @(26, 7, 7),
    @(27, 7, 7),
    @(28, 6, 6),
    @(29, 7, 7),
    @(30, 8, 8),
    @(31, 5, 6),
    @(32, 7, 7),
    @(33, 7, 7),
    @(34, 6, 6),
    @(35, 6, 6),
    @(36, 5, 5),
    @(37, 6, 6),
    @(38, 3, 3),
    @(39, 5, 5),
    @(40, 6, 6),
    @(41, 4, 4),
    @(42, 6, 6),
    @(43, 5, 5),
    @(44, 4, 4)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
